# reg_test_input_alt.xlsx: add "Qf" (meq/ml) as an alternative to "Qm"
# (meq/g) on the "params" sheet, and make "params" (cell C2) the active
# sheet/selection instead of "ions" (cell E4).

$wb = $excel.ActiveWorkbook
$paramsSheet = $wb.Worksheets.Item("params")

# Replace the "Qm" / meq/g row with the new "Qf" / meq/ml row.
$paramsSheet.Range("A2").Value = "Qf"
$paramsSheet.Range("B2").Value = 0.58589999999999998
$paramsSheet.Range("C2").Value = "meq/ml"

# "params" becomes the active sheet, with C2 selected.
$paramsSheet.Activate()
$paramsSheet.Range("C2").Select()
